# "THE BD IS REAL"
# The sheet originally had a single column A with 23 rows:
#   A1 = "Texte" (header, special style)
#   A2..A23 = question text (shared strings), each wrapped/vertically centered
#
# The edit:
#   1. Removes the "Texte" header row (row 1) entirely.
#   2. Inserts a brand-new column A that numbers the remaining 22 questions 1..22.
#   3. The original question text now lives in column B (same wrap/vertical-center
#      formatting it already had).
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the old header row ("Texte"). This shifts every other row up by one,
#    carrying their row heights/styles along with them.
$ws.Rows(1).Delete()

# 2. Make room for the new numbering column; old column A (now holding the
#    question text) becomes column B, keeping its formatting/width.
$ws.Columns("A").Insert()

# 3. Populate the new column A with plain row numbers 1..22.
for ($i = 1; $i -le 22; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
}

# Update the selection to mirror the final sheet state.
$ws.Range("A1:A22").Select() | Out-Null
